# Actualización automática 2025-10-24 10:30:09
#
# Applies the monthly sales-data correction for MEGAMAFERS S.A. (advisor
# HIDALGO HIDALGO PEDRO GUSTAVO): the "octubre" / PORCELANATO figure moves
# from 0 to -1581.47, and all of the dependent roll-up totals across the
# three sheets are refreshed to match.

$wb = $excel.ActiveWorkbook

$wsVentasPorGrupo   = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual     = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento     = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- VENTAS POR GRUPO: PORCELANATO sale for MEGAMAFERS S.A. (row 13) ---
$wsVentasPorGrupo.Range("M13").Value = -1581.47

# --- VENTA MENSUAL: octubre sale for MEGAMAFERS S.A. (row 13) ---
$wsVentaMensual.Range("F13").Value = -1581.47

# Column F (octubre) widens slightly to fit the new negative value.
$wsVentaMensual.Range("F1").ColumnWidth = 13.17

# Totals row (row 23) recalculated: SUM(F2:F22)
$wsVentaMensual.Range("F23").Value = 3519.52

# --- CUMPLIMIENTO MENSUAL: PORCELANATO group row (row 12) ---
$wsCumplimiento.Range("D12").Value = 2966.56
$wsCumplimiento.Range("E12").Value = 34773.18
$wsCumplimiento.Range("F12").Value = 0.078605734962668

# --- CUMPLIMIENTO MENSUAL: TOTAL row (row 14) ---
$wsCumplimiento.Range("D14").Value = 3519.52
$wsCumplimiento.Range("E14").Value = 51905.2214788039
$wsCumplimiento.Range("F14").Value = 0.06350088256786857
